$wb = $excel.ActiveWorkbook

# ALC row 132 (Leve Item ID 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 10961.75  # H132: 13953.667 -> 10961.75
$ws.Cells.Item(132, 9).Value = 2325  # I132: 2513.125 -> 2325
$ws.Cells.Item(132, 10).Value = 23916.875  # J132: 27028.572 -> 23916.875
$ws.Cells.Item(132, 11).Value = 6975  # K132: 7539.375 -> 6975
$ws.Cells.Item(132, 12).Value = 71750.625  # L132: 81085.716 -> 71750.625
$ws.Cells.Item(132, 13).Value = -4445  # M132: -5009.375 -> -4445
$ws.Cells.Item(132, 14).Value = -76810.625  # N132: -86145.716 -> -76810.625

# ARM row 2 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2455.3333  # H2: 2192.25 -> 2455.3333
$ws.Cells.Item(2, 9).Value = 3186.75  # I2: 2768.5557 -> 3186.75
$ws.Cells.Item(2, 10).Value = 992.5  # J2: 463.33334 -> 992.5
$ws.Cells.Item(2, 11).Value = 3186.75  # K2: 2768.5557 -> 3186.75
$ws.Cells.Item(2, 12).Value = 992.5  # L2: 463.33334 -> 992.5
$ws.Cells.Item(2, 13).Value = -3073.75  # M2: -2655.5557 -> -3073.75
$ws.Cells.Item(2, 14).Value = -1218.5  # N2: -689.33334 -> -1218.5

# ARM row 32 (Leve Item ID 44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3115.2754  # H32: 618.8099999999999 -> 3115.2754
$ws.Cells.Item(32, 9).Value = 3226.6775  # I32: 622.48456 -> 3226.6775
$ws.Cells.Item(32, 10).Value = 2128.5715  # J32: 500 -> 2128.5715
$ws.Cells.Item(32, 11).Value = 3226.6775  # K32: 622.48456 -> 3226.6775
$ws.Cells.Item(32, 12).Value = 2128.5715  # L32: 500 -> 2128.5715
$ws.Cells.Item(32, 13).Value = -2939.6775  # M32: -335.48456 -> -2939.6775
$ws.Cells.Item(32, 14).Value = -2702.5715  # N32: -1074 -> -2702.5715

# ARM row 74 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 18316884  # H74: 18825670 -> 18316884
$ws.Cells.Item(74, 9).Value = 21459368  # I74: 22889952 -> 21459368
$ws.Cells.Item(74, 11).Value = 21459368  # K74: 22889952 -> 21459368
$ws.Cells.Item(74, 13).Value = -21458494  # M74: -22889078 -> -21458494

# ARM row 77 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 18316884  # H77: 18825670 -> 18316884
$ws.Cells.Item(77, 9).Value = 21459368  # I77: 22889952 -> 21459368
$ws.Cells.Item(77, 11).Value = 107296840  # K77: 114449760 -> 107296840
$ws.Cells.Item(77, 13).Value = -107292472  # M77: -114445392 -> -107292472

# ARM row 116 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 2455.3333  # H116: 2192.25 -> 2455.3333
$ws.Cells.Item(116, 9).Value = 3186.75  # I116: 2768.5557 -> 3186.75
$ws.Cells.Item(116, 10).Value = 992.5  # J116: 463.33334 -> 992.5
$ws.Cells.Item(116, 11).Value = 3186.75  # K116: 2768.5557 -> 3186.75
$ws.Cells.Item(116, 12).Value = 992.5  # L116: 463.33334 -> 992.5
$ws.Cells.Item(116, 13).Value = -892.75  # M116: -474.5556999999999 -> -892.75
$ws.Cells.Item(116, 14).Value = -5580.5  # N116: -5051.33334 -> -5580.5

# ARM row 132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 2392.8823  # H132: 2343.6572 -> 2392.8823
$ws.Cells.Item(132, 9).Value = 1290.6666  # I132: 1251.875 -> 1290.6666
$ws.Cells.Item(132, 11).Value = 3871.9998  # K132: 3755.625 -> 3871.9998
$ws.Cells.Item(132, 13).Value = -1341.9998  # M132: -1225.625 -> -1341.9998

# BSM row 3 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2455.3333  # H3: 2192.25 -> 2455.3333
$ws.Cells.Item(3, 9).Value = 3186.75  # I3: 2768.5557 -> 3186.75
$ws.Cells.Item(3, 10).Value = 992.5  # J3: 463.33334 -> 992.5
$ws.Cells.Item(3, 11).Value = 3186.75  # K3: 2768.5557 -> 3186.75
$ws.Cells.Item(3, 12).Value = 992.5  # L3: 463.33334 -> 992.5
$ws.Cells.Item(3, 13).Value = -3072.75  # M3: -2654.5557 -> -3072.75
$ws.Cells.Item(3, 14).Value = -1220.5  # N3: -691.33334 -> -1220.5

# BSM row 80 (Leve Item ID 13747)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 294.125  # H80: 324.42856 -> 294.125
$ws.Cells.Item(80, 10).Value = 138.33333  # J80: 166.5 -> 138.33333
$ws.Cells.Item(80, 12).Value = 138.33333  # L80: 166.5 -> 138.33333
$ws.Cells.Item(80, 14).Value = -2134.33333  # N80: -2162.5 -> -2134.33333

# BSM row 83 (Leve Item ID 13747)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(83, 8).Value = 294.125  # H83: 324.42856 -> 294.125
$ws.Cells.Item(83, 10).Value = 138.33333  # J83: 166.5 -> 138.33333
$ws.Cells.Item(83, 12).Value = 691.6666499999999  # L83: 832.5 -> 691.6666499999999
$ws.Cells.Item(83, 14).Value = -10675.66665  # N83: -10816.5 -> -10675.66665

# BSM row 107 (Leve Item ID 27706)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 914.5  # H107: 1269.238 -> 914.5
$ws.Cells.Item(107, 9).Value = 883.2308  # I107: 1323.7646 -> 883.2308
$ws.Cells.Item(107, 10).Value = 1050  # J107: 1037.5 -> 1050
$ws.Cells.Item(107, 11).Value = 883.2308  # K107: 1323.7646 -> 883.2308
$ws.Cells.Item(107, 12).Value = 1050  # L107: 1037.5 -> 1050
$ws.Cells.Item(107, 13).Value = 1036.7692  # M107: 596.2354 -> 1036.7692
$ws.Cells.Item(107, 14).Value = -4890  # N107: -4877.5 -> -4890

# CRP row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3406.5686  # H31: 2477.1733 -> 3406.5686
$ws.Cells.Item(31, 9).Value = 1266.2727  # I31: 989.9394 -> 1266.2727
$ws.Cells.Item(31, 10).Value = 5030.241  # J31: 3645.7144 -> 5030.241
$ws.Cells.Item(31, 11).Value = 1266.2727  # K31: 989.9394 -> 1266.2727
$ws.Cells.Item(31, 12).Value = 5030.241  # L31: 3645.7144 -> 5030.241
$ws.Cells.Item(31, 13).Value = -971.2727  # M31: -694.9394 -> -971.2727
$ws.Cells.Item(31, 14).Value = -5620.241  # N31: -4235.7144 -> -5620.241

# CRP row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 3406.5686  # H34: 2477.1733 -> 3406.5686
$ws.Cells.Item(34, 9).Value = 1266.2727  # I34: 989.9394 -> 1266.2727
$ws.Cells.Item(34, 10).Value = 5030.241  # J34: 3645.7144 -> 5030.241
$ws.Cells.Item(34, 11).Value = 1266.2727  # K34: 989.9394 -> 1266.2727
$ws.Cells.Item(34, 12).Value = 5030.241  # L34: 3645.7144 -> 5030.241
$ws.Cells.Item(34, 13).Value = -1064.2727  # M34: -787.9394 -> -1064.2727
$ws.Cells.Item(34, 14).Value = -5434.241  # N34: -4049.7144 -> -5434.241

# CRP row 45 (Leve Item ID 2026)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(45, 8).Value = 30000  # H45: 3000 -> 30000
$ws.Cells.Item(45, 9).Value = 0  # I45: 3000 -> 0
$ws.Cells.Item(45, 10).Value = 30000  # J45: 0 -> 30000
$ws.Cells.Item(45, 11).Value = 0  # K45: 3000 -> 0
$ws.Cells.Item(45, 12).Value = 30000  # L45: 0 -> 30000
$ws.Cells.Item(45, 13).ClearContents()  # M45: -2407 -> (removed)
$ws.Cells.Item(45, 14).Value = -31186  # N45: None -> -31186

# CRP row 86 (Leve Item ID 12584)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 4180.1113  # H86: 4007.7896 -> 4180.1113
$ws.Cells.Item(86, 9).Value = 2370.1667  # I86: 2257.5386 -> 2370.1667
$ws.Cells.Item(86, 11).Value = 2370.1667  # K86: 2257.5386 -> 2370.1667
$ws.Cells.Item(86, 13).Value = -1247.1667  # M86: -1134.5386 -> -1247.1667

# CRP row 89 (Leve Item ID 12584)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 4180.1113  # H89: 4007.7896 -> 4180.1113
$ws.Cells.Item(89, 9).Value = 2370.1667  # I89: 2257.5386 -> 2370.1667
$ws.Cells.Item(89, 11).Value = 11850.8335  # K89: 11287.693 -> 11850.8335
$ws.Cells.Item(89, 13).Value = -6234.833500000001  # M89: -5671.692999999999 -> -6234.833500000001

# CRP row 99 (Leve Item ID 36198)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 52013.45  # H99: 37647.82 -> 52013.45
$ws.Cells.Item(99, 9).Value = 60699.94  # I99: 57362.168 -> 60699.94
$ws.Cells.Item(99, 10).Value = 2790  # J99: 2162 -> 2790
$ws.Cells.Item(99, 11).Value = 60699.94  # K99: 57362.168 -> 60699.94
$ws.Cells.Item(99, 12).Value = 2790  # L99: 2162 -> 2790
$ws.Cells.Item(99, 13).Value = -59201.94  # M99: -55864.168 -> -59201.94
$ws.Cells.Item(99, 14).Value = -5786  # N99: -5158 -> -5786

# CRP row 105 (Leve Item ID 19928)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 1212.762  # H105: 1198.5454 -> 1212.762
$ws.Cells.Item(105, 9).Value = 928.5333000000001  # I105: 926.75 -> 928.5333000000001
$ws.Cells.Item(105, 11).Value = 928.5333000000001  # K105: 926.75 -> 928.5333000000001
$ws.Cells.Item(105, 13).Value = 818.4666999999999  # M105: 820.25 -> 818.4666999999999

# CRP row 126 (Leve Item ID 36198)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 52013.45  # H126: 37647.82 -> 52013.45
$ws.Cells.Item(126, 9).Value = 60699.94  # I126: 57362.168 -> 60699.94
$ws.Cells.Item(126, 10).Value = 2790  # J126: 2162 -> 2790
$ws.Cells.Item(126, 11).Value = 182099.82  # K126: 172086.504 -> 182099.82
$ws.Cells.Item(126, 12).Value = 8370  # L126: 6486 -> 8370
$ws.Cells.Item(126, 13).Value = -179629.82  # M126: -169616.504 -> -179629.82
$ws.Cells.Item(126, 14).Value = -13310  # N126: -11426 -> -13310

# CRP row 132 (Leve Item ID 44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 26318284  # H132: 41669676 -> 26318284
$ws.Cells.Item(132, 9).Value = 29413846  # I132: 55557570 -> 29413846
$ws.Cells.Item(132, 10).Value = 5999  # J132: 5999.3335 -> 5999
$ws.Cells.Item(132, 11).Value = 88241538  # K132: 166672710 -> 88241538
$ws.Cells.Item(132, 12).Value = 17997  # L132: 17998.0005 -> 17997
$ws.Cells.Item(132, 13).Value = -88239008  # M132: -166670180 -> -88239008
$ws.Cells.Item(132, 14).Value = -23057  # N132: -23058.0005 -> -23057

# CRP row 134 (Leve Item ID 44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 31251866  # H134: 31251990 -> 31251866
$ws.Cells.Item(134, 9).Value = 45455556  # I134: 55556572 -> 45455556
$ws.Cells.Item(134, 10).Value = 3742.8  # J134: 3244.7144 -> 3742.8
$ws.Cells.Item(134, 11).Value = 136366668  # K134: 166669716 -> 136366668
$ws.Cells.Item(134, 12).Value = 11228.4  # L134: 9734.143199999999 -> 11228.4
$ws.Cells.Item(134, 13).Value = -136364133  # M134: -166667181 -> -136364133
$ws.Cells.Item(134, 14).Value = -16298.4  # N134: -14804.1432 -> -16298.4

# CUL row 12 (Leve Item ID 4854)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 31.731707  # H12: 35.054054 -> 31.731707
$ws.Cells.Item(12, 9).Value = 19.727272  # I12: 22.6 -> 19.727272
$ws.Cells.Item(12, 10).Value = 36.133335  # J12: 39.666668 -> 36.133335
$ws.Cells.Item(12, 11).Value = 59.181816  # K12: 67.80000000000001 -> 59.181816
$ws.Cells.Item(12, 12).Value = 108.400005  # L12: 119.000004 -> 108.400005
$ws.Cells.Item(12, 13).Value = 113.818184  # M12: 105.2 -> 113.818184
$ws.Cells.Item(12, 14).Value = -454.400005  # N12: -465.000004 -> -454.400005

# CUL row 113 (Leve Item ID 27843)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 655.9737  # H113: 687.32434 -> 655.9737
$ws.Cells.Item(113, 9).Value = 640.3929000000001  # I113: 668.88 -> 640.3929000000001
$ws.Cells.Item(113, 10).Value = 699.6  # J113: 725.75 -> 699.6
$ws.Cells.Item(113, 11).Value = 1921.1787  # K113: 2006.64 -> 1921.1787
$ws.Cells.Item(113, 12).Value = 2098.8  # L113: 2177.25 -> 2098.8
$ws.Cells.Item(113, 13).Value = 248.8212999999998  # M113: 163.3600000000001 -> 248.8212999999998
$ws.Cells.Item(113, 14).Value = -6438.8  # N113: -6517.25 -> -6438.8

# CUL row 119 (Leve Item ID 27873)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(119, 8).Value = 2462  # H119: 340 -> 2462
$ws.Cells.Item(119, 9).Value = 677  # I119: 340 -> 677
$ws.Cells.Item(119, 10).Value = 6032  # J119: 0 -> 6032
$ws.Cells.Item(119, 11).Value = 2031  # K119: 1020 -> 2031
$ws.Cells.Item(119, 12).Value = 18096  # L119: 0 -> 18096
$ws.Cells.Item(119, 13).Value = 2807  # M119: 3818 -> 2807
$ws.Cells.Item(119, 14).Value = -27772  # N119: None -> -27772

# CUL row 123 (Leve Item ID 36037)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(123, 8).Value = 3308.5715  # H123: 9486 -> 3308.5715
$ws.Cells.Item(123, 9).Value = 1830  # I123: 30030 -> 1830
$ws.Cells.Item(123, 10).Value = 3900  # J123: 4350 -> 3900
$ws.Cells.Item(123, 11).Value = 5490  # K123: 90090 -> 5490
$ws.Cells.Item(123, 12).Value = 11700  # L123: 13050 -> 11700
$ws.Cells.Item(123, 13).Value = -3040  # M123: -87640 -> -3040
$ws.Cells.Item(123, 14).Value = -16600  # N123: -17950 -> -16600

# CUL row 126 (Leve Item ID 36045)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(126, 8).Value = 4739.2  # H126: 7379.3 -> 4739.2
$ws.Cells.Item(126, 9).Value = 680  # I126: 8837.143 -> 680
$ws.Cells.Item(126, 10).Value = 8798.4  # J126: 3977.6667 -> 8798.4
$ws.Cells.Item(126, 11).Value = 2040  # K126: 26511.429 -> 2040
$ws.Cells.Item(126, 12).Value = 26395.2  # L126: 11933.0001 -> 26395.2
$ws.Cells.Item(126, 13).Value = 2900  # M126: -21571.429 -> 2900
$ws.Cells.Item(126, 14).Value = -36275.2  # N126: -21813.0001 -> -36275.2

# CUL row 131 (Leve Item ID 36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1097.2667  # H131: 1111.2667 -> 1097.2667
$ws.Cells.Item(131, 9).Value = 752.2941  # I131: 795.625 -> 752.2941
$ws.Cells.Item(131, 10).Value = 1233.6511  # J131: 1226.0454 -> 1233.6511
$ws.Cells.Item(131, 11).Value = 2256.8823  # K131: 2386.875 -> 2256.8823
$ws.Cells.Item(131, 12).Value = 3700.9533  # L131: 3678.1362 -> 3700.9533
$ws.Cells.Item(131, 13).Value = 2783.1177  # M131: 2653.125 -> 2783.1177
$ws.Cells.Item(131, 14).Value = -13780.9533  # N131: -13758.1362 -> -13780.9533

# CUL row 132 (Leve Item ID 43972)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 2081.7  # H132: 1730.7693 -> 2081.7
$ws.Cells.Item(132, 9).Value = 2127.2307  # I132: 2292.5 -> 2127.2307
$ws.Cells.Item(132, 10).Value = 1997.1428  # J132: 1249.2858 -> 1997.1428
$ws.Cells.Item(132, 11).Value = 19145.0763  # K132: 20632.5 -> 19145.0763
$ws.Cells.Item(132, 12).Value = 17974.2852  # L132: 11243.5722 -> 17974.2852
$ws.Cells.Item(132, 13).Value = -16615.0763  # M132: -18102.5 -> -16615.0763
$ws.Cells.Item(132, 14).Value = -23034.2852  # N132: -16303.5722 -> -23034.2852

# CUL row 133 (Leve Item ID 44073)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(133, 8).Value = 3616.4167  # H133: 4628.5713 -> 3616.4167
$ws.Cells.Item(133, 9).Value = 2613.375  # I133: 5010 -> 2613.375
$ws.Cells.Item(133, 10).Value = 5622.5  # J133: 4342.5 -> 5622.5
$ws.Cells.Item(133, 11).Value = 7840.125  # K133: 15030 -> 7840.125
$ws.Cells.Item(133, 12).Value = 16867.5  # L133: 13027.5 -> 16867.5
$ws.Cells.Item(133, 13).Value = -2780.125  # M133: -9970 -> -2780.125
$ws.Cells.Item(133, 14).Value = -26987.5  # N133: -23147.5 -> -26987.5

# CUL row 134 (Leve Item ID 44074)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(134, 8).Value = 2927.8572  # H134: 4000.3845 -> 2927.8572
$ws.Cells.Item(134, 9).Value = 1482.4546  # I134: 3010.7273 -> 1482.4546
$ws.Cells.Item(134, 10).Value = 8227.666999999999  # J134: 9443.5 -> 8227.666999999999
$ws.Cells.Item(134, 11).Value = 4447.3638  # K134: 9032.1819 -> 4447.3638
$ws.Cells.Item(134, 12).Value = 24683.001  # L134: 28330.5 -> 24683.001
$ws.Cells.Item(134, 13).Value = 622.6361999999999  # M134: -3962.1819 -> 622.6361999999999
$ws.Cells.Item(134, 14).Value = -34823.001  # N134: -38470.5 -> -34823.001

# CUL row 137 (Leve Item ID 44088)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value = 11515.083  # H137: 10650.692 -> 11515.083
$ws.Cells.Item(137, 9).Value = 1491.8  # I137: 1326.5 -> 1491.8
$ws.Cells.Item(137, 10).Value = 18674.572  # J137: 18642.857 -> 18674.572
$ws.Cells.Item(137, 11).Value = 4475.4  # K137: 3979.5 -> 4475.4
$ws.Cells.Item(137, 12).Value = 56023.716  # L137: 55928.571 -> 56023.716
$ws.Cells.Item(137, 13).Value = 624.6000000000004  # M137: 1120.5 -> 624.6000000000004
$ws.Cells.Item(137, 14).Value = -66223.716  # N137: -66128.571 -> -66223.716

# GSM row 122 (Leve Item ID 36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1075.25  # H122: 1014.375 -> 1075.25
$ws.Cells.Item(122, 9).Value = 1026.75  # I122: 935.6667 -> 1026.75
$ws.Cells.Item(122, 10).Value = 1123.75  # J122: 1061.6 -> 1123.75
$ws.Cells.Item(122, 11).Value = 3080.25  # K122: 2807.0001 -> 3080.25
$ws.Cells.Item(122, 12).Value = 3371.25  # L122: 3184.8 -> 3371.25
$ws.Cells.Item(122, 13).Value = -630.25  # M122: -357.0001000000002 -> -630.25
$ws.Cells.Item(122, 14).Value = -8271.25  # N122: -8084.799999999999 -> -8271.25

# GSM row 132 (Leve Item ID 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 4315135  # H132: 4469207 -> 4315135
$ws.Cells.Item(132, 9).Value = 7817326.5  # I132: 8338407.5 -> 7817326.5
$ws.Cells.Item(132, 11).Value = 23451979.5  # K132: 25015222.5 -> 23451979.5
$ws.Cells.Item(132, 13).Value = -23449449.5  # M132: -25012692.5 -> -23449449.5

# LTW row 51 (Leve Item ID 3423)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(51, 8).Value = 21950  # H51: 0 -> 21950
$ws.Cells.Item(51, 10).Value = 21950  # J51: 0 -> 21950
$ws.Cells.Item(51, 12).Value = 21950  # L51: 0 -> 21950
$ws.Cells.Item(51, 14).Value = -22906  # N51: None -> -22906

# LTW row 53 (Leve Item ID 3866)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(53, 8).Value = 18600  # H53: 12566.667 -> 18600
$ws.Cells.Item(53, 9).Value = 0  # I53: 500 -> 0
$ws.Cells.Item(53, 11).Value = 0  # K53: 500 -> 0
$ws.Cells.Item(53, 13).ClearContents()  # M53: 18 -> (removed)

# LTW row 68 (Leve Item ID 12563)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 2520.4614  # H68: 2362.0334 -> 2520.4614
$ws.Cells.Item(68, 9).Value = 2330.4  # I68: 2154.8333 -> 2330.4
$ws.Cells.Item(68, 10).Value = 2779.6365  # J68: 2672.8333 -> 2779.6365
$ws.Cells.Item(68, 11).Value = 2330.4  # K68: 2154.8333 -> 2330.4
$ws.Cells.Item(68, 12).Value = 2779.6365  # L68: 2672.8333 -> 2779.6365
$ws.Cells.Item(68, 13).Value = -1581.4  # M68: -1405.8333 -> -1581.4
$ws.Cells.Item(68, 14).Value = -4277.636500000001  # N68: -4170.8333 -> -4277.636500000001

# LTW row 71 (Leve Item ID 12563)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 2520.4614  # H71: 2362.0334 -> 2520.4614
$ws.Cells.Item(71, 9).Value = 2330.4  # I71: 2154.8333 -> 2330.4
$ws.Cells.Item(71, 10).Value = 2779.6365  # J71: 2672.8333 -> 2779.6365
$ws.Cells.Item(71, 11).Value = 11652  # K71: 10774.1665 -> 11652
$ws.Cells.Item(71, 12).Value = 13898.1825  # L71: 13364.1665 -> 13898.1825
$ws.Cells.Item(71, 13).Value = -7908  # M71: -7030.166499999999 -> -7908
$ws.Cells.Item(71, 14).Value = -21386.1825  # N71: -20852.1665 -> -21386.1825

# WVR row 50 (Leve Item ID 3421)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(50, 8).Value = 19542  # H50: 19694.666 -> 19542
$ws.Cells.Item(50, 10).Value = 19542  # J50: 19694.666 -> 19542
$ws.Cells.Item(50, 12).Value = 19542  # L50: 19694.666 -> 19542
$ws.Cells.Item(50, 14).Value = -20804  # N50: -20956.666 -> -20804

# WVR row 81 (Leve Item ID 12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2283  # H81: 2584.6365 -> 2283
$ws.Cells.Item(81, 9).Value = 829.1111  # I81: 990.1429000000001 -> 829.1111
$ws.Cells.Item(81, 10).Value = 4900  # J81: 5375 -> 4900
$ws.Cells.Item(81, 11).Value = 1658.2222  # K81: 1980.2858 -> 1658.2222
$ws.Cells.Item(81, 12).Value = 9800  # L81: 10750 -> 9800
$ws.Cells.Item(81, 13).Value = -597.2221999999999  # M81: -919.2858000000001 -> -597.2221999999999
$ws.Cells.Item(81, 14).Value = -11922  # N81: -12872 -> -11922

# WVR row 84 (Leve Item ID 12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(84, 8).Value = 2283  # H84: 2584.6365 -> 2283
$ws.Cells.Item(84, 9).Value = 829.1111  # I84: 990.1429000000001 -> 829.1111
$ws.Cells.Item(84, 10).Value = 4900  # J84: 5375 -> 4900
$ws.Cells.Item(84, 11).Value = 8291.110999999999  # K84: 9901.429 -> 8291.110999999999
$ws.Cells.Item(84, 12).Value = 49000  # L84: 53750 -> 49000
$ws.Cells.Item(84, 13).Value = -2987.110999999999  # M84: -4597.429 -> -2987.110999999999
$ws.Cells.Item(84, 14).Value = -59608  # N84: -64358 -> -59608

# WVR row 107 (Leve Item ID 27746)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 861  # H107: 886.3333 -> 861
$ws.Cells.Item(107, 9).Value = 797.75  # I107: 800 -> 797.75
$ws.Cells.Item(107, 10).Value = 911.6  # J107: 929.5 -> 911.6
$ws.Cells.Item(107, 11).Value = 2393.25  # K107: 2400 -> 2393.25
$ws.Cells.Item(107, 12).Value = 2734.8  # L107: 2788.5 -> 2734.8
$ws.Cells.Item(107, 13).Value = -473.25  # M107: -480 -> -473.25
$ws.Cells.Item(107, 14).Value = -6574.8  # N107: -6628.5 -> -6574.8

# WVR row 122 (Leve Item ID 36208)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 93057.82000000001  # H122: 102278.4 -> 93057.82000000001
$ws.Cells.Item(122, 9).Value = 126791.375  # I122: 144782.72 -> 126791.375
$ws.Cells.Item(122, 11).Value = 380374.125  # K122: 434348.16 -> 380374.125
$ws.Cells.Item(122, 13).Value = -377924.125  # M122: -431898.16 -> -377924.125

# WVR row 126 (Leve Item ID 36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1125.0555  # H126: 1203.091 -> 1125.0555
$ws.Cells.Item(126, 9).Value = 988.6774  # I126: 1046.5172 -> 988.6774
$ws.Cells.Item(126, 10).Value = 1970.6  # J126: 2338.25 -> 1970.6
$ws.Cells.Item(126, 11).Value = 2966.0322  # K126: 3139.5516 -> 2966.0322
$ws.Cells.Item(126, 12).Value = 5911.799999999999  # L126: 7014.75 -> 5911.799999999999
$ws.Cells.Item(126, 13).Value = -496.0322000000001  # M126: -669.5515999999998 -> -496.0322000000001
$ws.Cells.Item(126, 14).Value = -10851.8  # N126: -11954.75 -> -10851.8
